$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (row 1 title) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 15:52"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1011408
$ws.Range("C4").Value = 1052
$ws.Range("D4").Value = 139418
$ws.Range("E4").Value = 815135
$ws.Range("F4").Value = 14187
$ws.Range("G4").Value = 58
$ws.Range("H4").Value = 56855

# --- Row 42: Serbia ---
$ws.Range("F42").Value = 79

# --- Row 72: Islandia ---
$ws.Range("B72").Value = 1795
$ws.Range("C72").Value = 3
$ws.Range("E72").Value = 161

# --- Rows 148/149: Aruba and Sierra Leona swap places (Sierra Leona's
#     case count overtakes Aruba's), plus updated data for both ---
$ws.Range("A148").Value = "Sierra Leona"
$ws.Range("B148").Value = 104
$ws.Range("C148").Value = 5
$ws.Range("D148").Value = 12
$ws.Range("E148").Value = 88
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 4

$ws.Range("A149").Value = "Aruba"
$ws.Range("B149").Value = 100
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 73
$ws.Range("E149").Value = 25
$ws.Range("F149").Value = 4
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 2
